{"js": "// LCWA log: bump the firmware version mentioned in the two existing log\n// rows from 9.02.01 -> 9.02.02, and fill in a new log entry (Time,\n// Box/CW, Change) in the first blank row of the table. The Date cell of\n// that row is intentionally left blank, matching the target edit.\n\n// 1) Correct \"9.02.01\" -> \"9.02.02\" everywhere it appears in the body\n// (both existing log rows reference the same version string).\nconst versionHits = context.document.body.search(\"9.02.01\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\n\nversionHits.items.forEach((hit) => {\n  hit.insertText(\"9.02.02\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// 2) Locate the log table and the first fully empty row (a blank entry\n// pre-provisioned in the table) to record the new log line.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nlet targetRow = -1;\nfor (let i = 1; i < table.values.length; i++) {\n  const row = table.values[i];\n  if (row.every((cell) => cell === \"\")) {\n    targetRow = i;\n    break;\n  }\n}\n\nif (targetRow === -1) {\n  throw new Error(\"No blank log row found to fill in.\");\n}\n\n// Columns: 0 = Date, 1 = Time, 2 = Box/CW, 3 = Change.\n// Insert at the end of each (currently empty) cell paragraph so the\n// existing run/paragraph formatting is preserved instead of being\n// replaced by a freshly-created default run.\nconst timeCell = table.getCell(targetRow, 1);\ntimeCell.body.paragraphs.getFirst().insertText(\"06:08\", Word.InsertLocation.end);\n\nconst boxCell = table.getCell(targetRow, 2);\nboxCell.body.paragraphs.getFirst().insertText(\"LC24\", Word.InsertLocation.end);\n\nconst changeCell = table.getCell(targetRow, 3);\nchangeCell.body.paragraphs.getFirst().insertText(\"speedserver = 9686\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# LCWA log: bump the firmware version mentioned in the two existing log\n# rows from 9.02.01 -> 9.02.02, and fill in a new log entry (Time,\n# Box/CW, Change) in the first blank row of the table. The Date cell of\n# that row is intentionally left blank, matching the target edit.\n\n$d = $word.ActiveDocument\n\n# 1) Correct \"9.02.01\" -> \"9.02.02\" everywhere it appears in the document.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"9.02.01\"\n$find.Replacement.Text = \"9.02.02\"\n$find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n\n# 2) Locate the log table and the first fully empty data row (a blank\n# entry pre-provisioned in the table) to record the new log line.\n$table = $d.Tables.Item(1)\n\n$targetRow = -1\nfor ($r = 2; $r -le $table.Rows.Count; $r++) {\n    $isBlank = $true\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cellText = $table.Cell($r, $c).Range.Text -replace \"[\\x07\\x0d]\", \"\"\n        if ($cellText -ne \"\") {\n            $isBlank = $false\n            break\n        }\n    }\n    if ($isBlank) {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -eq -1) {\n    throw \"No blank log row found to fill in.\"\n}\n\n# Columns: 1 = Date, 2 = Time, 3 = Box/CW, 4 = Change.\n# Assigning directly to Range.Text keeps the existing (empty) run/\n# paragraph formatting instead of inserting a freshly-formatted run.\n$table.Cell($targetRow, 2).Range.Text = \"06:08\"\n$table.Cell($targetRow, 3).Range.Text = \"LC24\"\n$table.Cell($targetRow, 4).Range.Text = \"speedserver = 9686\"\n"}
